$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at 228, pushing existing rows 228-336 down to 229-337.
$ws.Rows(228).Insert()

# Populate the newly inserted row 228 with the new record.
$ws.Range("A228").Value = 5
$ws.Range("B228").Value = "Macroferia Regional de Talca"
$ws.Range("C228").Value = "Maule"
$ws.Range("D228").Value = 45029
$ws.Range("E228").Value = 7
$ws.Range("F228").Value = 100112024
$ws.Range("G228").Value = "Choclo"
$ws.Range("H228").Value = "Choclero"
$ws.Range("I228").Value = "Primera"
$ws.Range("J228").Value = 20000
$ws.Range("K228").Value = 300
$ws.Range("L228").Value = 300
$ws.Range("M228").Value = 300
$ws.Range("N228").Value = "$/unidad"
$ws.Range("O228").Value = "Región del Maule"
$ws.Range("P228").Value = 300
$ws.Range("Q228").Value = 1
$ws.Range("R228").Value = "Hortaliza"
